$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.613.48"
$ws.Range("E2").Value = "  -2.28%  "
$ws.Range("D3").Value = "'1.842.31"
$ws.Range("E3").Value = "  -1.47%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'314.73"
$ws.Range("E5").Value = "  -1.23%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").Value = "'0.4263"
$ws.Range("E7").Value = "  -2.73%  "
$ws.Range("D8").Value = "'0.3641"
$ws.Range("E8").Value = "  -1.53%  "
$ws.Range("D9").Value = "'45.67"
$ws.Range("E9").Value = "  +0.81%  "
$ws.Range("D10").Value = "'0.07280"
$ws.Range("E10").Value = "  -3.06%  "
$ws.Range("D11").Value = "'0.8954"
$ws.Range("D12").Value = "'20.61"
$ws.Range("E12").Value = "  -3.83%  "
$ws.Range("D13").Value = "'1.864.46"
$ws.Range("E13").Value = "  -2.57%  "
$ws.Range("E14").Value = "  -1.15%  "
$ws.Range("D15").Value = "'6.558"
$ws.Range("E15").Value = "  -2.47%  "
$ws.Range("D16").Value = "'0.06887"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "'78.26"
$ws.Range("E18").Value = "  -4.61%  "
$ws.Range("D19").Value = "'0.000008878"
$ws.Range("E19").Value = "  -2.01%  "
$ws.Range("D20").Value = "'1.002"
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("E21").Value = "  -2.37%  "
$ws.Range("D22").Value = "'27.655.77"
$ws.Range("E22").Value = "  -2.01%  "
$ws.Range("D23").Value = "'4.980"
$ws.Range("E23").Value = "  -2.92%  "
$ws.Range("D24").Value = "'10.54"
$ws.Range("E24").Value = "  -2.25%  "
$ws.Range("D25").Value = "'2.096.93"
$ws.Range("E25").Value = "  -0.92%  "
$ws.Range("D26").Value = "'2.041"
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("D27").Value = "'154.58"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").Value = "'18.31"
$ws.Range("E28").Value = "  -0.60%  "
$ws.Range("D29").Value = "'5.232"
$ws.Range("E29").Value = "  -1.49%  "
$ws.Range("D30").Value = "'1.860"
$ws.Range("E30").Value = "  +7.66%  "
$ws.Range("D31").Value = "'111.72"
$ws.Range("E31").Value = "  -1.79%  "
$ws.Range("D32").Value = "'0.08892"
$ws.Range("E32").Value = "  -1.82%  "
$ws.Range("D33").Value = "'0.7764"
$ws.Range("E33").Value = "  -2.85%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'4.575"
$ws.Range("E34").Value = "  -5.47%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'2.990"
$ws.Range("E35").Value = "  +2.17%  "
$ws.Range("D36").Value = "'1.100"
$ws.Range("E36").Value = "  -6.11%  "
$ws.Range("D37").Value = "'1.000"
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("D38").Value = "'0.05437"
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").Value = "'1.099"
$ws.Range("E39").Value = "  -2.06%  "
$ws.Range("D40").Value = "'0.01929"
$ws.Range("E40").Value = "  -1.19%  "
$ws.Range("D41").Value = "'2.817"
$ws.Range("E41").Value = "  -3.24%  "
$ws.Range("D42").Value = "'0.5075"
$ws.Range("E42").Value = "  -3.42%  "
$ws.Range("D43").Value = "'6.804"
$ws.Range("E43").Value = "  -4.04%  "
$ws.Range("D44").Value = "'0.1643"
$ws.Range("E44").Value = "  -2.15%  "
$ws.Range("D45").Value = "'8.233"
$ws.Range("E45").Value = "  -5.82%  "
$ws.Range("D46").Value = "'0.06652"
$ws.Range("E46").Value = "  -1.43%  "
$ws.Range("D47").Value = "'10.41"
$ws.Range("E47").Value = "  -1.28%  "
$ws.Range("D48").Value = "'106.16"
$ws.Range("E48").Value = "  -1.66%  "
$ws.Range("E49").Value = "  -3.39%  "
$ws.Range("D50").Value = "'0.9995"
$ws.Range("E50").Value = "  -0.26%  "
$ws.Range("E51").Value = "  -2.53%  "
